{"js": "// Replace each two-digit multiplication equation in the document body\n// with its updated value, per the commit diff.\nconst replacements = [\n  [\"41\u00d793=3813\", \"63\u00d718=1134\"],\n  [\"79\u00d746=3634\", \"51\u00d773=3723\"],\n  [\"16\u00d795=1520\", \"62\u00d720=1240\"],\n  [\"73\u00d711=803\", \"38\u00d788=3344\"],\n  [\"74\u00d785=6290\", \"57\u00d762=3534\"],\n  [\"21\u00d742=882\", \"83\u00d750=4150\"],\n  [\"37\u00d747=1739\", \"30\u00d779=2370\"],\n  [\"55\u00d732=1760\", \"58\u00d720=1160\"],\n  [\"11\u00d727=297\", \"69\u00d714=966\"],\n  [\"17\u00d760=1020\", \"45\u00d723=1035\"],\n  [\"93\u00d786=7998\", \"21\u00d768=1428\"],\n  [\"71\u00d736=2556\", \"45\u00d789=4005\"],\n  [\"59\u00d712=708\", \"97\u00d757=5529\"],\n  [\"68\u00d711=748\", \"97\u00d764=6208\"],\n  [\"38\u00d745=1710\", \"29\u00d793=2697\"],\n  [\"19\u00d776=1444\", \"72\u00d719=1368\"],\n  [\"17\u00d732=544\", \"13\u00d747=611\"],\n  [\"27\u00d744=1188\", \"88\u00d797=8536\"],\n  [\"36\u00d735=1260\", \"59\u00d779=4661\"],\n  [\"33\u00d775=2475\", \"57\u00d755=3135\"],\n  [\"99\u00d750=4950\", \"95\u00d793=8835\"],\n  [\"17\u00d776=1292\", \"94\u00d774=6956\"],\n  [\"16\u00d759=944\", \"33\u00d737=1221\"],\n  [\"37\u00d766=2442\", \"68\u00d777=5236\"],\n  [\"99\u00d787=8613\", \"17\u00d718=306\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication equation in the document with its\n# updated value, per the commit diff.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"41\u00d793=3813\", \"63\u00d718=1134\"),\n    @(\"79\u00d746=3634\", \"51\u00d773=3723\"),\n    @(\"16\u00d795=1520\", \"62\u00d720=1240\"),\n    @(\"73\u00d711=803\", \"38\u00d788=3344\"),\n    @(\"74\u00d785=6290\", \"57\u00d762=3534\"),\n    @(\"21\u00d742=882\", \"83\u00d750=4150\"),\n    @(\"37\u00d747=1739\", \"30\u00d779=2370\"),\n    @(\"55\u00d732=1760\", \"58\u00d720=1160\"),\n    @(\"11\u00d727=297\", \"69\u00d714=966\"),\n    @(\"17\u00d760=1020\", \"45\u00d723=1035\"),\n    @(\"93\u00d786=7998\", \"21\u00d768=1428\"),\n    @(\"71\u00d736=2556\", \"45\u00d789=4005\"),\n    @(\"59\u00d712=708\", \"97\u00d757=5529\"),\n    @(\"68\u00d711=748\", \"97\u00d764=6208\"),\n    @(\"38\u00d745=1710\", \"29\u00d793=2697\"),\n    @(\"19\u00d776=1444\", \"72\u00d719=1368\"),\n    @(\"17\u00d732=544\", \"13\u00d747=611\"),\n    @(\"27\u00d744=1188\", \"88\u00d797=8536\"),\n    @(\"36\u00d735=1260\", \"59\u00d779=4661\"),\n    @(\"33\u00d775=2475\", \"57\u00d755=3135\"),\n    @(\"99\u00d750=4950\", \"95\u00d793=8835\"),\n    @(\"17\u00d776=1292\", \"94\u00d774=6956\"),\n    @(\"16\u00d759=944\", \"33\u00d737=1221\"),\n    @(\"37\u00d766=2442\", \"68\u00d777=5236\"),\n    @(\"99\u00d787=8613\", \"17\u00d718=306\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
